# Update evaluation results on the "QuantitativeMetrics" sheet of the
# test code generation evaluation workbook (UC3.4.1_TC1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Compilation success note: replace the previous remark with the new one.
$ws.Range("C5").Value = "Calling a not existing method"

# Code BLEU score updated value.
$ws.Range("B12").Value = 0.3036414372992564

# Code BLEU breakdown note updated to match the new score / dataflow match.
$ws.Range("C12").Value = "{'codebleu': 0.3036414372992564, 'ngram_match_score': 0.14347354720083563, 'weighted_ngram_match_score': 0.15355275219867576, 'syntax_match_score': 0.551948051948052, 'dataflow_match_score': 0.3655913978494624}"
